$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1) Finish off the last row of the "Week 14" section (row 120): the day of
#    44541 (12/11/2021) gets logged as 0.75 hrs of JS101, "Finish 3 small
#    problems".
# ---------------------------------------------------------------------------
$ws.Range("B120").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C120").Value = 0.75
$ws.Range("D120").Value = "Finish 3 small problems"

# ---------------------------------------------------------------------------
# 2) Insert 9 fresh rows right before the grand-total row (old row 122) to
#    hold a new "Week 15" section: a header row, 7 date rows and a weekly
#    total row.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 9; $i++) {
    $ws.Rows.Item(122).Insert()
}

# Row 122: "Week 15" section header - match the look of the other section
# headers (e.g. row 113, "Week 14"): bold left-aligned date-formatted cell
# in column A, plain cells in C/E, nothing in B/D.
$ws.Range("A113").Copy()
$ws.Range("A122").PasteSpecial(-4122)
$ws.Range("C113").Copy()
$ws.Range("C122").PasteSpecial(-4122)
$ws.Range("E113").Copy()
$ws.Range("E122").PasteSpecial(-4122)
$ws.Range("B122").Clear()
$ws.Range("D122").Clear()
$ws.Range("A122").Value = "Week 15"

# Rows 123-129: the 7 days of Week 15 (12/12/2021 - 12/18/2021), still empty
# of hours/notes, just like the freshly-started week used to look before it
# had any entries (B/D cleared out, A holds the date).
$dates = @(44542, 44543, 44544, 44545, 44546, 44547, 44548)
$r = 123
foreach ($d in $dates) {
    $ws.Cells.Item($r, 1).Value = $d
    $ws.Range("B$r").Clear()
    $ws.Range("D$r").Clear()
    $r++
}

# Row 130: Weekly Total row for Week 15, mirroring the other weekly-total
# rows (e.g. row 121): bold "Weekly Total" label + SUM formula over the
# week's Hours column.
$ws.Range("B130").Value = "Weekly Total"
$ws.Range("D130").Formula = "=SUM(C123:C129)"

# ---------------------------------------------------------------------------
# 3) Resize Table1 so the new rows (and the shifted-down totals row) are part
#    of the table, matching the new A1:E131 extent.
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range("A1:E131"))

# ---------------------------------------------------------------------------
# 4) Move the active selection down onto the new week's first Course cell,
#    mirroring where editing left off.
# ---------------------------------------------------------------------------
[void]$ws.Range("B123").Select()
